# Applies the "fix a mistake when counting adds in polygonise" edit:
# updates the raw "adds" counts (column G) in both tables on sheet "baseline",
# and the dependent literal copies of the computed ratios in columns K and L
# that feed the scatter chart. Formulas in column H recompute automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline")

# --- Table 1 (rows 5-14): columns F (total ticks), G (adds), H (=G/F) ---
$G1 = @(30776, 161308, 431100, 911708, 1625212, 2624508, 3968828, 5711356, 7853724, 10524284)

for ($i = 0; $i -lt $G1.Length; $i++) {
    $row = 5 + $i
    $ws.Range("G$row").Value = $G1[$i]
}

# Column K holds a literal (non-formula) copy of H5:H14 that is plotted on
# the chart as the "radius=0.4" series' Y values.
$K1 = @(0.37262829870136749, 0.31310500495510679, 0.27687994643519476, 0.26435552127627188, 0.2526449643156391, 0.24316747345456863, 0.23392777572138224, 0.22746174253110946, 0.2218500674725864, 0.21755588962102432)

for ($i = 0; $i -lt $K1.Length; $i++) {
    $row = 5 + $i
    $ws.Range("K$row").Value = $K1[$i]
}

# Column L holds a literal (non-formula) copy of H21:H30 (table 2's ratios),
# plotted as the "radius=0.05" series' Y values.
$L1 = @(0.20320460489974385, 0.18805316785198462, 0.19450935912462144, 0.19384699442068654, 0.19504045908765708, 0.19451522506218943, 0.19462093043898493, 0.19391802149488666, 0.19442028596481808, 0.19320328245856411)

for ($i = 0; $i -lt $L1.Length; $i++) {
    $row = 5 + $i
    $ws.Range("L$row").Value = $L1[$i]
}

# --- Table 2 (rows 21-30): columns F (total ticks), G (adds), H (=G/F) ---
$G2 = @(8504, 64540, 220536, 516572, 1010584, 1735452, 2762648, 4118716, 5866776, 8032508)

for ($i = 0; $i -lt $G2.Length; $i++) {
    $row = 21 + $i
    $ws.Range("G$row").Value = $G2[$i]
}

# Recalculate so the H column formulas (=G/F) and the chart caches refresh.
$excel.CalculateFullRebuild()

# --- Chart cosmetic tweak: drop the light-grey outline around the whole
# chart object (chart area border -> no line). ---
$chartObj = $ws.ChartObjects(1)
$chartObj.Chart.ChartArea.Format.Line.Visible = $false

# --- Restore the view to where the author left it when saving. ---
$window = $excel.ActiveWindow
$ws.Activate()
$ws.Range("D13").Select()
$window.ScrollRow = 13
$window.ScrollColumn = 4
$ws.Range("I10").Select()
